$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 25000112
$ws.Range("I39").Value = 50000024
$ws.Range("J39").Value = 200
$ws.Range("K39").Value = 150000072
$ws.Range("L39").Value = 600
$ws.Range("M39").Value = -149999776
$ws.Range("N39").Value = -1192

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 6666.6665
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 6666.6665
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 3033880.2
$ws.Range("I112").Value = 2100
$ws.Range("J112").Value = 3337058.2
$ws.Range("K112").Value = 6300
$ws.Range("L112").Value = 10011174.6
$ws.Range("M112").Value = -5192
$ws.Range("N112").Value = -10013390.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 150000
$ws.Range("I116").Value = 150000
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 150000
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -146558
$ws.Range("N116").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 1422.0769
$ws.Range("I125").Value = 1309.25
$ws.Range("J125").Value = 1472.2222
$ws.Range("K125").Value = 11783.25
$ws.Range("L125").Value = 13249.9998
$ws.Range("M125").Value = -9323.25
$ws.Range("N125").Value = -18169.9998

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H126").Value = 141998.5
$ws.Range("J126").Value = 141998.5
$ws.Range("L126").Value = 141998.5
$ws.Range("N126").Value = -151878.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 126666.664
$ws.Range("J133").Value = 126666.664
$ws.Range("L133").Value = 126666.664
$ws.Range("N133").Value = -136786.664

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2967.75
$ws.Range("I74").Value = 2434.7144
$ws.Range("K74").Value = 2434.7144
$ws.Range("M74").Value = -1560.7144

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2967.75
$ws.Range("I77").Value = 2434.7144
$ws.Range("K77").Value = 12173.572
$ws.Range("M77").Value = -7805.572

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 4502.86
$ws.Range("I122").Value = 3391.6
$ws.Range("K122").Value = 10174.8
$ws.Range("M122").Value = -7724.799999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 22490512
$ws.Range("I134").Value = 5105490
$ws.Range("K134").Value = 15316470
$ws.Range("M134").Value = -15313935

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2675
$ws.Range("I16").Value = 2450.3333
$ws.Range("J16").Value = 2809.8
$ws.Range("K16").Value = 2450.3333
$ws.Range("L16").Value = 2809.8
$ws.Range("M16").Value = -2163.3333
$ws.Range("N16").Value = -3383.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 145930.14
$ws.Range("I107").Value = 501505.5
$ws.Range("K107").Value = 501505.5
$ws.Range("M107").Value = -499585.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 2675
$ws.Range("I113").Value = 2450.3333
$ws.Range("J113").Value = 2809.8
$ws.Range("K113").Value = 2450.3333
$ws.Range("L113").Value = 2809.8
$ws.Range("M113").Value = -280.3332999999998
$ws.Range("N113").Value = -7149.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H125").Value = 112666.336
$ws.Range("I125").Value = 99999
$ws.Range("J125").Value = 119000
$ws.Range("K125").Value = 99999
$ws.Range("L125").Value = 119000
$ws.Range("M125").Value = -97539
$ws.Range("N125").Value = -123920

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1927.84
$ws.Range("I132").Value = 1864.8
$ws.Range("K132").Value = 5594.4
$ws.Range("M132").Value = -3064.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2046.7028
$ws.Range("I134").Value = 1522.1613
$ws.Range("J134").Value = 4756.8335
$ws.Range("K134").Value = 4566.4839
$ws.Range("L134").Value = 14270.5005
$ws.Range("M134").Value = -2031.4839
$ws.Range("N134").Value = -19340.5005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 57227250
$ws.Range("I4").Value = 35267652
$ws.Range("K4").Value = 105802956
$ws.Range("M4").Value = -105802844

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 3999.75
$ws.Range("J63").Value = 4999.6665
$ws.Range("L63").Value = 14998.9995
$ws.Range("N63").Value = -16496.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H66").Value = 3999.75
$ws.Range("J66").Value = 4999.6665
$ws.Range("L66").Value = 44996.9985
$ws.Range("N66").Value = -52484.9985

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1981.5
$ws.Range("I129").Value = 1964
$ws.Range("J129").Value = 1999
$ws.Range("K129").Value = 5892
$ws.Range("L129").Value = 5997
$ws.Range("M129").Value = -892
$ws.Range("N129").Value = -15997

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 26627442
$ws.Range("I138").Value = 2307.8
$ws.Range("K138").Value = 6923.400000000001
$ws.Range("M138").Value = -1783.400000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 945.6667
$ws.Range("I107").Value = 906.6667
$ws.Range("J107").Value = 1023.6667
$ws.Range("K107").Value = 906.6667
$ws.Range("L107").Value = 1023.6667
$ws.Range("M107").Value = 1013.3333
$ws.Range("N107").Value = -4863.6667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H124").Value = 74999.5
$ws.Range("J124").Value = 74999.5
$ws.Range("L124").Value = 74999.5
$ws.Range("N124").Value = -84819.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3799.7727
$ws.Range("I132").Value = 3837.8572
$ws.Range("K132").Value = 11513.5716
$ws.Range("M132").Value = -8983.571599999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H139").Value = 88823.5
$ws.Range("J139").Value = 88823.5
$ws.Range("L139").Value = 88823.5
$ws.Range("N139").Value = -99103.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 14000.667
$ws.Range("I40").Value = 14000.667
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 14000.667
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -13864.667
$ws.Range("N40").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 1000000000
$ws.Range("J43").Value = 1000000000
$ws.Range("L43").Value = 1000000000
$ws.Range("N43").Value = -1000000386

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3599.5
$ws.Range("J46").Value = 4012
$ws.Range("L46").Value = 4012
$ws.Range("N46").Value = -4388

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 634.44446
$ws.Range("I55").Value = 523.5
$ws.Range("K55").Value = 523.5
$ws.Range("M55").Value = -350.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3805.3845
$ws.Range("I122").Value = 3831.75
$ws.Range("J122").Value = 3489
$ws.Range("K122").Value = 11495.25
$ws.Range("L122").Value = 10467
$ws.Range("M122").Value = -9045.25
$ws.Range("N122").Value = -15367

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3610
$ws.Range("I132").Value = 2581.6667
$ws.Range("J132").Value = 5666.6665
$ws.Range("K132").Value = 7745.000100000001
$ws.Range("L132").Value = 16999.9995
$ws.Range("M132").Value = -5215.000100000001
$ws.Range("N132").Value = -22059.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 78427.5
$ws.Range("J46").Value = 78427.5
$ws.Range("L46").Value = 78427.5
$ws.Range("N46").Value = -78889.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3211.2173
$ws.Range("I122").Value = 2738.3684
$ws.Range("J122").Value = 5457.25
$ws.Range("K122").Value = 8215.1052
$ws.Range("L122").Value = 16371.75
$ws.Range("M122").Value = -5765.1052
$ws.Range("N122").Value = -21271.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H125").Value = 72961.60000000001
$ws.Range("J125").Value = 72961.60000000001
$ws.Range("L125").Value = 72961.60000000001
$ws.Range("N125").Value = -82801.60000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H134").Value = 78427.5
$ws.Range("J134").Value = 78427.5
$ws.Range("L134").Value = 235282.5
$ws.Range("N134").Value = -240352.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 20207.387
$ws.Range("I136").Value = 1698.2285
$ws.Range("J136").Value = 49653.773
$ws.Range("K136").Value = 5094.6855
$ws.Range("L136").Value = 148961.319
$ws.Range("M136").Value = -2544.6855
$ws.Range("N136").Value = -154061.319
